$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column from 2023-10-22 (45221)
# to 2023-10-25 (45224) for rows 2 through 5.
foreach ($r in 2..5) {
    $ws.Cells.Item($r, 3).Value = 45224
}
